$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Fix F2 / F7, which were previously stored as text-typed shared
# strings ("1905755" / "2260463") but should be plain numbers, matching
# the numeric F3:F6 cells already on the sheet. ---
$ws.Range("F2").Value = 1905755
$ws.Range("F7").Value = 2260463

# --- Extend the table with the common legislator/property metadata
# columns (G:M), matching the other sheets in this workbook. Copy the
# existing header/data styles (s="1" header, s="2" data) across instead
# of re-deriving new style indices. ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("G2:M7").PasteSpecial(-4122) | Out-Null

# Header row
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Data rows 2-7
$rows = 2..7
$indexValues = @(77, 78, 79, 80, 81, 82)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"
    $ws.Range("I$r").Value = "2011-12-21"
    $ws.Range("J$r").Value = "林明溱"
    $ws.Range("K$r").Value = 1706
    $ws.Range("L$r").Value = "tmp56941"
    $ws.Range("M$r").Value = $indexValues[$i]
}
